$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 6
$ws.Range("D2").Value = 4
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = -3
$ws.Range("G2").Value = -3
$ws.Range("H2").Value = 43

$ws.Range("B3").Value = 6
$ws.Range("C3").Value = 5
$ws.Range("D3").Value = 5
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = -1
$ws.Range("G3").Value = -5
$ws.Range("H3").Value = 21

$ws.Range("B4").Value = 6
$ws.Range("C4").Value = 8
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 7

$ws.Range("B5").Value = 5
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 6
$ws.Range("F5").Value = -4
$ws.Range("G5").Value = -2
$ws.Range("H5").Value = 54

$ws.Range("B6").Value = 8
$ws.Range("C6").Value = 9
$ws.Range("D6").Value = 6
$ws.Range("E6").Value = 5
$ws.Range("F6").Value = -2
$ws.Range("G6").Value = -4
$ws.Range("H6").Value = 32

$ws.Range("I1").Select() | Out-Null
